# SCRUM Planung C, E und J
#
# Sprint planning pass:
#  - pull the next three Backlog items into the Sprint-Backlog with
#    concrete time estimates / due dates / done-state
#  - mark the first two sprint items as finished ("Done") and add the
#    team's due date
#  - note the sprint's Focus Faktor
#  - flesh out the short Backlog task descriptions with the detailed
#    sprint-planning notes

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Backlog"
$ws2 = $wb.Worksheets.Item(2)   # "Sprint-Backlog"

# ---------------------------------------------------------------------
# 1) Sprint-Backlog: bring in Backlog task #3 ("Implementierung der
#    TicTacToe-Darstellung") as sprint row 5, already finished.
# ---------------------------------------------------------------------
$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = $ws1.Range("B5").Value2
$ws2.Range("C5").Value = $ws1.Range("C5").Value2
$ws2.Range("D5").Value = "10min"
$ws2.Range("E5").Value = "Done"

# ---------------------------------------------------------------------
# 2) Sprint-Backlog: mark the first two sprint items as "Done" as well.
# ---------------------------------------------------------------------
$ws2.Range("E3").Value = "Done"
$ws2.Range("E4").Value = "Done"

# ---------------------------------------------------------------------
# 3) Backlog: flesh out the short "Umsetzung ..." descriptions with the
#    detailed sprint-planning text.
# ---------------------------------------------------------------------
$ws1.Range("C7").Value = "Umsetzung der Rules im Framework. Standard TicToeMit. Mit Ausblick auf Zeitbegrenzung: wer schneller ist bekommt einen Bonus. Mit Ausblick: Auswahl ob Kreuz oder Kreis."
$ws1.Range("C6").Value = "Umsetzung des Painters im Framework. Standard TicTacToe. Kreuze, Kreise einfarbiger Hintergrund. Mit Ausblick auf skalierte Anpassung auf Bildschirmgröße. Mit Ausblick: Auswahl ob Kreuz oder Kreis."
$ws1.Range("C9").Value = "Umsetzung eines PC gesteuerten Players im Framework. Standard TicTacToe-Spieler. Ausblick auf ""Schwierigkeitsstufen""."
$ws1.Range("C8").Value = "Umsetzung eines menschlichen Players im Framework. Standard TicTacToe mit Mauswahl."

# ---------------------------------------------------------------------
# 4) Sprint-Backlog: note this sprint's Focus Faktor.
# ---------------------------------------------------------------------
$ws2.Range("D1").Value = "Focus Faktor: 0,5"
$ws2.Range("D1").HorizontalAlignment = -4108
$ws2.Range("D1").VerticalAlignment = -4108
$ws2.Range("D1").WrapText = $false

# ---------------------------------------------------------------------
# 5) Sprint-Backlog: bring in Backlog tasks #4 and #5 ("Implementierung
#    der TicTacToe-Regeln" / "... eines menschlichen TicTacToe-Spielers")
#    as sprint rows 6 and 7, with their estimated effort.
# ---------------------------------------------------------------------
$ws2.Range("A6").Value = 2
$ws2.Range("B6").Value = $ws1.Range("B6").Value2
$ws2.Range("C6").Value = $ws1.Range("C6").Value2
$ws2.Range("D6").Value = "250min"

$ws2.Range("A7").Value = 2
$ws2.Range("B7").Value = $ws1.Range("B8").Value2
$ws2.Range("C7").Value = $ws1.Range("C8").Value2
$ws2.Range("D7").Value = "100min"

# ---------------------------------------------------------------------
# 6) Due dates for the first two sprint items.
# ---------------------------------------------------------------------
$ws2.Range("D3").Value = "Due: 08.11.2021"
$ws2.Range("D4").Value = "Due: 08.11.2021"

# ---------------------------------------------------------------------
# 7) Cosmetics: center the new "Done" column (E) and size it like the
#    other narrow columns. Rows 6/7 don't carry a "Done" value (not yet
#    finished) but still get the column touched so the row's cell range
#    is recorded.
# ---------------------------------------------------------------------
$ws2.Range("E3:E5").HorizontalAlignment = -4108
$ws2.Range("E3:E5").VerticalAlignment = -4108
$ws2.Range("E3:E5").WrapText = $false
$ws2.Range("E6").Style = "Standard"
$ws2.Range("E7").Style = "Standard"
$ws2.Columns.Item(5).ColumnWidth = 8.88671875

# ---------------------------------------------------------------------
# 8) Update the selections / scroll position shown when the file is
#    reopened.
# ---------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("A8:XFD8").Select()

$ws2.Activate()
$ws2.Range("D7").Select()
